# Update NATMI LR-pair data (Jag1-Notch2) with newly recomputed TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 11.190867
$ws.Range("H2").Value = 33.572601
$ws.Range("I2").Value = 0.1514016037116739
$ws.Range("J2").Value = 0.1514016037116739
$ws.Range("M2").Value = 1.492477333333333
$ws.Range("N2").Value = 4.477432
$ws.Range("O2").Value = 0.02769484181536182
$ws.Range("P2").Value = 0.02769484181536182
$ws.Range("Q2").Value = 16.702115337848
$ws.Range("R2").Value = 150.319038040632
$ws.Range("S2").Value = 0.004193043465386905
$ws.Range("T2").Value = 0.004193043465386905

$ws.Range("G3").Value = 11.190867
$ws.Range("H3").Value = 33.572601
$ws.Range("I3").Value = 0.1514016037116739
$ws.Range("J3").Value = 0.1514016037116739
$ws.Range("O3").Value = 0.6282762845978157
$ws.Range("P3").Value = 0.6282762845978156
$ws.Range("Q3").Value = 378.898823085776
$ws.Range("R3").Value = 3410.089407771984
$ws.Range("S3").Value = 0.09512203706212131
$ws.Range("T3").Value = 0.0951220370621213

$ws.Range("G4").Value = 11.190867
$ws.Range("H4").Value = 33.572601
$ws.Range("I4").Value = 0.1514016037116739
$ws.Range("J4").Value = 0.1514016037116739
$ws.Range("N4").Value = 55.619234
$ws.Range("O4").Value = 0.3440288735868225
$ws.Range("P4").Value = 0.3440288735868225
$ws.Range("Q4").Value = 207.475816778626
$ws.Range("R4").Value = 1867.282351007634
$ws.Range("S4").Value = 0.05208652318416565
$ws.Range("T4").Value = 0.05208652318416565

$ws.Range("I5").Value = 0.2043341870182926
$ws.Range("J5").Value = 0.2043341870182926
$ws.Range("M5").Value = 1.492477333333333
$ws.Range("N5").Value = 4.477432
$ws.Range("O5").Value = 0.02769484181536182
$ws.Range("P5").Value = 0.02769484181536182
$ws.Range("Q5").Value = 22.54145976910667
$ws.Range("R5").Value = 202.87313792196
$ws.Range("S5").Value = 0.005659002986942174
$ws.Range("T5").Value = 0.005659002986942173

$ws.Range("I6").Value = 0.2043341870182926
$ws.Range("J6").Value = 0.2043341870182926
$ws.Range("O6").Value = 0.6282762845978157
$ws.Range("P6").Value = 0.6282762845978156
$ws.Range("S6").Value = 0.1283783238361681
$ws.Range("T6").Value = 0.1283783238361681

$ws.Range("I7").Value = 0.2043341870182926
$ws.Range("J7").Value = 0.2043341870182926
$ws.Range("N7").Value = 55.619234
$ws.Range("O7").Value = 0.3440288735868225
$ws.Range("P7").Value = 0.3440288735868225
$ws.Range("S7").Value = 0.07029686019518235
$ws.Range("T7").Value = 0.07029686019518235

$ws.Range("G8").Value = 47.62086333333333
$ws.Range("I8").Value = 0.6442642092700336
$ws.Range("J8").Value = 0.6442642092700336
$ws.Range("M8").Value = 1.492477333333333
$ws.Range("N8").Value = 4.477432
$ws.Range("O8").Value = 0.02769484181536182
$ws.Range("P8").Value = 0.02769484181536182
$ws.Range("Q8").Value = 71.07305911876443
$ws.Range("R8").Value = 639.65753206888
$ws.Range("S8").Value = 0.01784279536303275
$ws.Range("T8").Value = 0.01784279536303275

$ws.Range("G9").Value = 47.62086333333333
$ws.Range("I9").Value = 0.6442642092700336
$ws.Range("J9").Value = 0.6442642092700336
$ws.Range("O9").Value = 0.6282762845978157
$ws.Range("P9").Value = 0.6282762845978156
$ws.Range("S9").Value = 0.4047759236995263
$ws.Range("T9").Value = 0.4047759236995263

$ws.Range("G10").Value = 47.62086333333333
$ws.Range("I10").Value = 0.6442642092700336
$ws.Range("J10").Value = 0.6442642092700336
$ws.Range("N10").Value = 55.619234
$ws.Range("O10").Value = 0.3440288735868225
$ws.Range("P10").Value = 0.3440288735868225
$ws.Range("Q10").Value = 882.8786470062288
$ws.Range("R10").Value = 7945.907823056059
$ws.Range("S10").Value = 0.2216454902074746
$ws.Range("T10").Value = 0.2216454902074746

